# Update "Keamanan" workbook: daftar lalu lintas
$wb = $excel.ActiveWorkbook

# --- DaftarLaluLintas_Input sheet ---
$wsInput = $wb.Worksheets.Item("DaftarLaluLintas_Input")

# Update the row counter in A2 from 3 to 5
$wsInput.Range("A2").Value = 5

# Move the selection to B13 and make this the active (selected) sheet/tab
$wsInput.Activate()
$wsInput.Range("B13").Select()

# --- DaftarLaluLintas_Index sheet ---
$wsIndex = $wb.Worksheets.Item("DaftarLaluLintas_Index")

# Move the selection to C3 (this sheet is no longer the active tab)
$wsIndex.Range("C3").Select()

# Re-activate the input sheet so it remains the workbook's active tab
$wsInput.Activate()
